# Updates cryptos list with refreshed prices/volumes, and restores the
# original row order for Uniswap/InternetComputer and Monero/PEPE.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, [string]$value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = '@'
    $cell.Value = $value
    $cell.Style = 'Normal'
}

$ws.Range('D2').Value = '64.093.21'
$ws.Range('E2').Value = '  +2.68%  '

$ws.Range('D3').Value = '3.061.99'
$ws.Range('E3').Value = '  +2.04%  '

$ws.Range('E4').Value = '  +0.10%  '

Set-CellText 'D5' '559.42'
$ws.Range('E5').Value = '  +2.69%  '

Set-CellText 'D6' '143.85'
$ws.Range('E6').Value = '  +3.34%  '

Set-CellText 'D7' '1.00'

$ws.Range('D8').Value = '3.064.23'
$ws.Range('E8').Value = '  +2.29%  '

Set-CellText 'D9' '0.512'
$ws.Range('E9').Value = '  +4.76%  '

Set-CellText 'D10' '0.156'
$ws.Range('E10').Value = '  +5.37%  '

Set-CellText 'D11' '6.10'
$ws.Range('E11').Value = '  -9.92%  '

$ws.Range('E12').Value = '  +7.83%  '

$ws.Range('E13').Value = '  +5.14%  '

Set-CellText 'D14' '35.75'
$ws.Range('E14').Value = '  +5.08%  '

$ws.Range('D15').Value = '3.563.30'
$ws.Range('E15').Value = '  +2.36%  '

$ws.Range('D16').Value = '64.194.78'
$ws.Range('E16').Value = '  +2.80%  '

$ws.Range('D17').Value = '3.068.69'
$ws.Range('E17').Value = '  +2.17%  '

$ws.Range('E18').Value = '  +2.12%  '

$ws.Range('E19').Value = '  +3.25%  '

Set-CellText 'D20' '479.09'
$ws.Range('E20').Value = '  +2.01%  '

$ws.Range('E21').Value = '  +4.55%  '

Set-CellText 'D22' '0.681'
$ws.Range('E22').Value = '  +4.38%  '

$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-CellText 'D23' '7.63'
$ws.Range('E23').Value = '  +6.12%  '

$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-CellText 'D24' '14.38'
$ws.Range('E24').Value = '  +14.00%  '

Set-CellText 'D25' '82.20'
$ws.Range('E25').Value = '  +3.27%  '

$ws.Range('E26').Value = '  +0.14%  '

$ws.Range('E27').Value = '  +2.94%  '

Set-CellText 'D28' '8.01'
$ws.Range('E28').Value = '  +5.18%  '

Set-CellText 'D29' '2.05'
$ws.Range('E29').Value = '  +1.65%  '

Set-CellText 'D30' '1.00'
$ws.Range('E30').Value = '  -0.06%  '

Set-CellText 'D31' '26.39'
$ws.Range('E31').Value = '  +3.75%  '

$ws.Range('E32').Value = '  +1.05%  '

Set-CellText 'D33' '2.44'
$ws.Range('E33').Value = '  +3.97%  '

$ws.Range('E34').Value = '  +2.55%  '

$ws.Range('E35').Value = '  +7.12%  '

Set-CellText 'D36' '54.95'
$ws.Range('E36').Value = '  +0.58%  '

$ws.Range('E37').Value = '  +4.99%  '

Set-CellText 'D38' '446.46'
$ws.Range('E38').Value = '  -0.81%  '

$ws.Range('E39').Value = '  +0.44%  '

Set-CellText 'D40' '2.87'
$ws.Range('E40').Value = '  +11.09%  '

$ws.Range('D41').Value = '2.994.28'
$ws.Range('E41').Value = '  +1.29%  '

Set-CellText 'D42' '8.26'
$ws.Range('E42').Value = '  +2.52%  '

Set-CellText 'D43' '0.116'
$ws.Range('E43').Value = '  +1.42%  '

$ws.Range('E44').Value = '  +4.74%  '

Set-CellText 'D45' '0.262'
$ws.Range('E45').Value = '  +6.06%  '

Set-CellText 'D46' '2.16'
$ws.Range('E46').Value = '  +7.86%  '

$ws.Range('E48').Value = '  +4.08%  '

$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-CellText 'D49' '119.21'
$ws.Range('E49').Value = '  +3.40%  '

$ws.Range('B50').Value = 'PEPE'
$ws.Range('C50').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D50').Value = '0.0₃0520'
$ws.Range('E50').Value = '  +4.86%  '

$ws.Range('E51').Value = '  +3.56%  '
